$d = $word.ActiveDocument

# --- 1. Locate the real "AUTOEVALUACIÓN" / "BIBLIOGRAFÍA" Heading-1 paragraphs ---
# (not the Table-of-Contents lines, which share the same words but live in
#  "TOC 1" styled paragraphs)
$autoEvalIdx = -1
$biblioIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $sn = $p.Style.NameLocal
    $txt = $p.Range.Text
    if ($sn -eq "Heading 1" -and $txt -like "*AUTOEVALUACI*") {
        $autoEvalIdx = $i
    }
    if ($sn -eq "Heading 1" -and $txt -like "*BIBLIOGRAF*") {
        $biblioIdx = $i
    }
}

$pStart = $d.Paragraphs.Item($autoEvalIdx)
$pEnd = $d.Paragraphs.Item($biblioIdx)
$full = $d.Range($pStart.Range.Start, $pEnd.Range.End)

# --- 2. Replace the 3-paragraph block (AUTOEVALUACIÓN / blank / BIBLIOGRAFÍA)
#        with the new 6-paragraph block: the "Dificultades del desarrollo"
#        sub-heading + its body text, a blank line, the (moved) AUTOEVALUACIÓN
#        heading, a blank line, and the (moved) BIBLIOGRAFÍA heading
#        (dropping the stray <w:lastRenderedPageBreak/>). ---
$body = ""
$body += '<w:p><w:pPr><w:pStyle w:val="textogeneral"/><w:rPr><w:rFonts w:ascii="Lucida Sans" w:hAnsi="Lucida Sans"/><w:b/><w:bCs/><w:color w:val="23B8DC"/><w:sz w:val="24"/><w:szCs w:val="48"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Lucida Sans" w:hAnsi="Lucida Sans"/><w:b/><w:bCs/><w:color w:val="23B8DC"/><w:sz w:val="24"/><w:szCs w:val="48"/></w:rPr><w:t>Dificultades del desarrollo</w:t></w:r></w:p>'
$body += '<w:p><w:pPr><w:pStyle w:val="textogeneral"/><w:rPr><w:rFonts w:ascii="Lucida Sans Unicode" w:eastAsia="Times New Roman" w:hAnsi="Lucida Sans Unicode" w:cs="Lucida Sans Unicode"/><w:kern w:val="0"/><w:sz w:val="24"/><w:lang w:eastAsia="es-ES" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Lucida Sans Unicode" w:eastAsia="Times New Roman" w:hAnsi="Lucida Sans Unicode" w:cs="Lucida Sans Unicode"/><w:kern w:val="0"/><w:sz w:val="24"/><w:lang w:eastAsia="es-ES" w:bidi="ar-SA"/></w:rPr><w:t>Al comprobar el tipo de variables boleanas he tenido que cambiar la manera de comprobar</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Lucida Sans Unicode" w:eastAsia="Times New Roman" w:hAnsi="Lucida Sans Unicode" w:cs="Lucida Sans Unicode"/><w:kern w:val="0"/><w:sz w:val="24"/><w:lang w:eastAsia="es-ES" w:bidi="ar-SA"/></w:rPr><w:t>.</w:t></w:r></w:p>'
$body += '<w:p><w:pPr><w:pStyle w:val="textogeneral"/></w:pPr></w:p>'
$body += '<w:p><w:pPr><w:pStyle w:val="Ttulo1"/></w:pPr><w:r><w:t>AUTOEVALUACIÓN</w:t></w:r></w:p>'
$body += '<w:p><w:pPr><w:pStyle w:val="textogeneral"/></w:pPr></w:p>'
$body += '<w:p><w:pPr><w:pStyle w:val="Ttulo1"/></w:pPr><w:r><w:t>BIBLIOGRAFÍA</w:t></w:r></w:p>'

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body>' + $body + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$full.InsertXML($xml)

# --- 3. Cached PAGE field in the footer of the last section drops from 2 to 1 ---
foreach ($sec in $d.Sections) {
    $ftr = $sec.Footers.Item(1)
    [void]$ftr.Range.Find.Execute("2", $true, $false, $false, $false, $false, $true, 1, $false, "1", 2)
}
